$wb = $excel.ActiveWorkbook
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcn.Range("D4").Value = "2016-02-25 05:27:18"
$dede.Range("D4").Value = "2016-02-25 05:27:30"
